$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E40").Value = 2
$ws.Range("E44").Value = 1
$ws.Range("E60").Value = 2
$excel.Calculate()

$cos = $ws.ChartObjects()
$co = $cos.Item(2)
$chart = $co.Chart
$series = $chart.SeriesCollection(2)
$series.Values = "=Hoja1!`$E`$67:`$E`$71"
